# Daily attendance processing - swap the order of "Recorded By" values in
# column G from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$updated = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $updated = $updated + 1
    }
}

Write-Host "Updated $updated cell(s) in column G."
